{"js": "// Apply the diff: update the date line and the 25 division-problem table\n// cells. Each entry is [oldText, newText]; every oldText string occurs\n// exactly once in the document, so an exact (case-sensitive, non-whole-word)\n// search-and-replace on each one reproduces the edit unambiguously.\nconst replacements = [\n  [\"2026-02-13 Friday\", \"2026-02-14 Saturday\"],\n  [\"294\u00f72=147, 0\", \"444\u00f73=148, 0\"],\n  [\"112\u00f76=18, 4\", \"621\u00f78=77, 5\"],\n  [\"833\u00f73=277, 2\", \"842\u00f78=105, 2\"],\n  [\"646\u00f73=215, 1\", \"917\u00f76=152, 5\"],\n  [\"156\u00f76=26, 0\", \"115\u00f79=12, 7\"],\n  [\"633\u00f77=90, 3\", \"556\u00f73=185, 1\"],\n  [\"840\u00f77=120, 0\", \"439\u00f78=54, 7\"],\n  [\"813\u00f74=203, 1\", \"906\u00f79=100, 6\"],\n  [\"208\u00f74=52, 0\", \"338\u00f79=37, 5\"],\n  [\"480\u00f77=68, 4\", \"797\u00f74=199, 1\"],\n  [\"496\u00f76=82, 4\", \"622\u00f74=155, 2\"],\n  [\"475\u00f74=118, 3\", \"665\u00f73=221, 2\"],\n  [\"579\u00f77=82, 5\", \"744\u00f75=148, 4\"],\n  [\"749\u00f72=374, 1\", \"520\u00f78=65, 0\"],\n  [\"592\u00f77=84, 4\", \"460\u00f77=65, 5\"],\n  [\"559\u00f74=139, 3\", \"382\u00f72=191, 0\"],\n  [\"183\u00f76=30, 3\", \"193\u00f76=32, 1\"],\n  [\"189\u00f74=47, 1\", \"700\u00f79=77, 7\"],\n  [\"608\u00f74=152, 0\", \"398\u00f76=66, 2\"],\n  [\"493\u00f73=164, 1\", \"542\u00f78=67, 6\"],\n  [\"794\u00f76=132, 2\", \"190\u00f73=63, 1\"],\n  [\"983\u00f78=122, 7\", \"627\u00f74=156, 3\"],\n  [\"137\u00f77=19, 4\", \"667\u00f73=222, 1\"],\n  [\"473\u00f73=157, 2\", \"375\u00f79=41, 6\"],\n  [\"716\u00f72=358, 0\", \"528\u00f75=105, 3\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$wdFindContinue = 1\n$wdReplaceAll = 2\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @('2026-02-13 Friday', '2026-02-14 Saturday'),\n    @('294\u00f72=147, 0', '444\u00f73=148, 0'),\n    @('112\u00f76=18, 4', '621\u00f78=77, 5'),\n    @('833\u00f73=277, 2', '842\u00f78=105, 2'),\n    @('646\u00f73=215, 1', '917\u00f76=152, 5'),\n    @('156\u00f76=26, 0', '115\u00f79=12, 7'),\n    @('633\u00f77=90, 3', '556\u00f73=185, 1'),\n    @('840\u00f77=120, 0', '439\u00f78=54, 7'),\n    @('813\u00f74=203, 1', '906\u00f79=100, 6'),\n    @('208\u00f74=52, 0', '338\u00f79=37, 5'),\n    @('480\u00f77=68, 4', '797\u00f74=199, 1'),\n    @('496\u00f76=82, 4', '622\u00f74=155, 2'),\n    @('475\u00f74=118, 3', '665\u00f73=221, 2'),\n    @('579\u00f77=82, 5', '744\u00f75=148, 4'),\n    @('749\u00f72=374, 1', '520\u00f78=65, 0'),\n    @('592\u00f77=84, 4', '460\u00f77=65, 5'),\n    @('559\u00f74=139, 3', '382\u00f72=191, 0'),\n    @('183\u00f76=30, 3', '193\u00f76=32, 1'),\n    @('189\u00f74=47, 1', '700\u00f79=77, 7'),\n    @('608\u00f74=152, 0', '398\u00f76=66, 2'),\n    @('493\u00f73=164, 1', '542\u00f78=67, 6'),\n    @('794\u00f76=132, 2', '190\u00f73=63, 1'),\n    @('983\u00f78=122, 7', '627\u00f74=156, 3'),\n    @('137\u00f77=19, 4', '667\u00f73=222, 1'),\n    @('473\u00f73=157, 2', '375\u00f79=41, 6'),\n    @('716\u00f72=358, 0', '528\u00f75=105, 3')\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $ok = $find.Execute($oldText, $true, $true, $false, $false, $false, $true, $wdFindContinue, $false, $newText, $wdReplaceAll)\n    if (-not $ok) {\n        throw \"Text not found: $oldText\"\n    }\n}\n"}
